$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Shift the N:P "Map Construction" key/value block down by one row ---
# Old layout (rows 9-14):
#   N9  NPC Count        O9  0   P9  =(60*60)-1 (3599)
#   N10 Extra map count  O10 0   P10 "-"
#   N11 Map Doors        O11 0   P11 3600
#   N12 Map data         O12 1   P12 3600
#   N13 (empty)
#   N14 Note: See example below...
#
# New layout (rows 9-15):
#   N9  Map Flags         O9  1   P9  "-"
#   N10 NPC Count         O10 0   P10 =(60*60)-1 (3599)
#   N11 Extra map count   O11 0   P11 "-"
#   N12 Map Doors         O12 0   P12 3600
#   N13 Map data          O13 1   P13 3600
#   N14 (empty)
#   N15 Note: See example below...

$ws.Range("N9").Value = "Map Flags"
$ws.Range("O9").Value = 1
$ws.Range("P9").Value = "-"

$ws.Range("N10").Value = "NPC Count"
$ws.Range("O10").Value = 0
$ws.Range("P10").Formula = "=(60*60)-1"

$ws.Range("N11").Value = "Extra map count"
$ws.Range("O11").Value = 0
$ws.Range("P11").Value = "-"

$ws.Range("N12").Value = "Map Doors"
$ws.Range("O12").Value = 0
$ws.Range("P12").Value = 3600

$ws.Range("N13").Value = "Map data"
$ws.Range("O13").Value = 1
$ws.Range("P13").Value = 3600

$ws.Range("N14").ClearContents()

$ws.Range("N15").Value = "Note: See example below on how the document MUST be formatted."

# --- Move the picture (documentation screenshot) down/over slightly ---
$shp = $ws.Shapes.Item(1)
$shp.Left = 829.810546875
$shp.Top = 239.25

# --- Update the view: selection + visible scroll position ---
$ws.Range("P9").Select()
